$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Miedź
$ws2 = $wb.Worksheets.Item(2)   # Aluminium
$ws3 = $wb.Worksheets.Item(3)   # Mosiądz
$ws4 = $wb.Worksheets.Item(4)   # Stal

# ---------------------------------------------------------------------------
# Row 1 headers - written in the exact order needed so the shared-strings
# table gets populated with indices 4..18 in the same sequence as the target
# workbook.
# ---------------------------------------------------------------------------

# Sheet1 (Miedź) new headers
$ws1.Range("F1").Value = "Średnia prędkość fali"   # -> shared string 4
$ws1.Range("G1").Value = "błąd bezwzględny"        # -> shared string 5
$ws1.Range("M1").Value = "Sp"                      # -> shared string 6
$ws1.Range("N1").Value = "V"                       # -> shared string 7
$ws1.Range("K1").Value = "l [m]"                   # -> shared string 8
$ws1.Range("J1").Value = "d [m]"                   # -> shared string 9
$ws1.Range("I1").Value = "dw [m]"                  # -> shared string 10
$ws1.Range("L1").Value = "m [kg]"                  # -> shared string 11
$ws1.Range("O1").Value = "ro"                      # -> shared string 12

# Sheet2 (Aluminium) new headers
$ws2.Range("F1").Value = "Średnia prędkość fali"
$ws2.Range("G1").Value = "błąd bezwzględny"
$ws2.Range("I1").Value = "h [m]"                   # -> shared string 13
$ws2.Range("J1").Value = "d [m]"
$ws2.Range("K1").Value = "m [kg]"
$ws2.Range("L1").Value = "Pp"                      # -> shared string 14
$ws2.Range("M1").Value = "V"
$ws2.Range("N1").Value = "ro"

# Sheet4 (Stal) new headers (done before sheet3 so that the "a/b/c [m]"
# strings land at indices 15..17, matching the target workbook)
$ws4.Range("F1").Value = "Średnia prędkość fali"
$ws4.Range("G1").Value = "błąd bezwzględny"
$ws4.Range("I1").Value = "a [m]"                   # -> shared string 15
$ws4.Range("J1").Value = "b [m]"                   # -> shared string 16
$ws4.Range("K1").Value = "c [m]"                   # -> shared string 17
$ws4.Range("L1").Value = "m [kg]"
$ws4.Range("M1").Value = "V"
$ws4.Range("N1").Value = "ro"

# Sheet3 (Mosiądz) new headers (reuse existing strings only)
$ws3.Range("F1").Value = "Średnia prędkość fali"
$ws3.Range("G1").Value = "błąd bezwzględny"
$ws3.Range("I1").Value = "h [m]"
$ws3.Range("J1").Value = "d [m]"
$ws3.Range("K1").Value = "m [kg]"
$ws3.Range("L1").Value = "Pp"
$ws3.Range("M1").Value = "V"
$ws3.Range("N1").Value = "ro"

# Add "Young" last - > shared string 18
$ws1.Range("J5").Value = "Young"
$ws2.Range("J5").Value = "Young"
$ws3.Range("J5").Value = "Young"
$ws4.Range("J5").Value = "Young"

# ---------------------------------------------------------------------------
# Sheet1 (Miedź)
# ---------------------------------------------------------------------------
$ws1.Range("F2").Formula = "=AVERAGE(D3:D7)"
$ws1.Range("G2").Formula = "=ABS(D2-`$F`$2)"
$ws1.Range("G3:G7").Formula = "=ABS(D3-`$F`$2)"

$ws1.Range("I2").Value = 0.0152
$ws1.Range("J2").Value = 0.01795
$ws1.Range("K2").Value = 1.801
$ws1.Range("L2").Value = 0.761
$ws1.Range("M2").Formula = "=(PI()/4)*((J2^2)-(I2^2))"
$ws1.Range("N2").Formula = "=M2*K2"
$ws1.Range("O2").Formula = "=L2/N2"
$ws1.Range("J6").Formula = "=O2*F2^2"

# styling: D2 gets the "applyNumberFormat" style, K2 gets red centered font
$ws1.Range("D2").NumberFormat = "General"
$ws1.Range("K2").Font.Color = 255
$ws1.Range("K2").HorizontalAlignment = -4108
$ws1.Range("K2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Sheet2 (Aluminium)
# ---------------------------------------------------------------------------
$ws2.Range("F2").Formula = "=AVERAGE(D2:D6)"
$ws2.Range("G2").Formula = "=ABS(D2-`$F`$2)"
$ws2.Range("G3:G7").Formula = "=ABS(D3-`$F`$2)"

$ws2.Range("I2").Value = 0.339
$ws2.Range("J2").Value = 0.0049
$ws2.Range("K2").Value = 0.024
$ws2.Range("L2").Formula = "=(PI()/4)*J2^2"
$ws2.Range("M2").Formula = "=I2*L2"
$ws2.Range("N2").Formula = "=K2/M2"
$ws2.Range("J6").Formula = "=N2*F2^2"

$ws2.Range("D7").NumberFormat = "General"

# ---------------------------------------------------------------------------
# Sheet3 (Mosiądz)
# ---------------------------------------------------------------------------
$ws3.Range("F2").Formula = "=AVERAGE(D2:D6)"
$ws3.Range("G2").Formula = "=ABS(D2-`$F`$2)"
$ws3.Range("G3:G7").Formula = "=ABS(D3-`$F`$2)"

$ws3.Range("I2").Value = 0.311
$ws3.Range("J2").Value = 0.0059
$ws3.Range("K2").Value = 0.074
$ws3.Range("L2").Formula = "=(PI()/4)*J2^2"
$ws3.Range("M2").Formula = "=I2*L2"
$ws3.Range("N2").Formula = "=K2/M2"
$ws3.Range("J6").Formula = "=N2*F2^2"

$ws3.Range("D7").NumberFormat = "General"

# ---------------------------------------------------------------------------
# Sheet4 (Stal)
# ---------------------------------------------------------------------------
$ws4.Range("F2").Formula = "=SUM(D2:D7)/A7"
$ws4.Range("G2").Formula = "=ABS(D2-`$F`$2)"
$ws4.Range("G3:G7").Formula = "=ABS(D3-`$F`$2)"

$ws4.Range("I2").Value = 0.019
$ws4.Range("J2").Value = 0.0141
$ws4.Range("K2").Value = 0.0142
$ws4.Range("L2").Value = 0.03186
$ws4.Range("M2").Formula = "=I2*J2*K2"
$ws4.Range("N2").Formula = "=L2/M2"
$ws4.Range("J6").Formula = "=N2*F2^2"

# ---------------------------------------------------------------------------
# Column widths for the newly-added columns
# ---------------------------------------------------------------------------
$ws1.Columns.Item(6).ColumnWidth = 9.88    # F
$ws1.Columns.Item(13).ColumnWidth = 12.59  # M
$ws1.Columns.Item(14).ColumnWidth = 11.02  # N
$ws1.Columns.Item(15).ColumnWidth = 9.74   # O

$ws2.Columns.Item(10).ColumnWidth = 9.17   # J
$ws2.Columns.Item(12).ColumnWidth = 11.17  # L
$ws2.Columns.Item(13).ColumnWidth = 11.17  # M
$ws2.Columns.Item(14).ColumnWidth = 9.31   # N

$ws3.Columns.Item(10).ColumnWidth = 11.17  # J
$ws3.Columns.Item(12).ColumnWidth = 11.17  # L
$ws3.Columns.Item(13).ColumnWidth = 11.17  # M

$ws4.Columns.Item(7).ColumnWidth = 9.02    # G
$ws4.Columns.Item(10).ColumnWidth = 10.17  # J
$ws4.Columns.Item(13).ColumnWidth = 11.17  # M

# ---------------------------------------------------------------------------
# Selections / active sheet (workbookView activeTab + tabSelected move to
# Mosiądz / sheet3, as in the target workbook)
# ---------------------------------------------------------------------------
$ws1.Range("I12").Select()
$ws2.Activate()
$ws2.Range("I11").Select()
$ws4.Activate()
$ws4.Range("I15").Select()
$ws3.Activate()
$ws3.Range("I9").Select()
